$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting from H1 so that the
# existing style ("s=1") is reused instead of a new style being created.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new numeric data for columns I and J (rows 2-9).
$iValues = @(5, 6, 7, 7, 6, 7, 10, 5)
$jValues = @(5, 7, 8, 8, 7, 8, 10, 5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
